# The workbook tracks weekly wholesale prices for "Zanahoria" (carrot) sold by
# "Agricola del Norte S.A. de Arica". Two brand-new weekly price records (for
# the reporting date 2021-09-21 / serial 44460) are inserted at the top of the
# data block (which starts at row 124), pushing every existing record down by
# two rows. The sheet's used range grows from A1:R188 to A1:R190.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the data block; Excel shifts every
# existing row (124-188) down to (126-190) and carries formatting along.
$ws.Rows("124:125").Insert()

# New row 124: "Primera" quality, same date/volume context, new date & prices.
$ws.Range("A124").Value = 1
$ws.Range("B124").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C124").Value = "Arica y Parinacota"
$ws.Range("D124").Value = 44460
$ws.Range("E124").Value = 15
$ws.Range("F124").Value = 100114013
$ws.Range("G124").Value = "Zanahoria"
$ws.Range("H124").Value = "Sin especificar"
$ws.Range("I124").Value = "Primera"
$ws.Range("J124").Value = 80
$ws.Range("K124").Value = 7000
$ws.Range("L124").Value = 7500
$ws.Range("M124").Value = 7250
$ws.Range("N124").Value = "$/saco 25 kilos"
$ws.Range("O124").Value = "Valle de Camiña"
$ws.Range("P124").Value = 290
$ws.Range("Q124").Value = 25
$ws.Range("R124").Value = "Hortaliza"

# New row 125: "Segunda" quality, same new date, new volume & prices.
$ws.Range("A125").Value = 1
$ws.Range("B125").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C125").Value = "Arica y Parinacota"
$ws.Range("D125").Value = 44460
$ws.Range("E125").Value = 15
$ws.Range("F125").Value = 100114013
$ws.Range("G125").Value = "Zanahoria"
$ws.Range("H125").Value = "Sin especificar"
$ws.Range("I125").Value = "Segunda"
$ws.Range("J125").Value = 60
$ws.Range("K125").Value = 5000
$ws.Range("L125").Value = 6000
$ws.Range("M125").Value = 5500
$ws.Range("N125").Value = "$/saco 25 kilos"
$ws.Range("O125").Value = "Valle de Camiña"
$ws.Range("P125").Value = 220
$ws.Range("Q125").Value = 25
$ws.Range("R125").Value = "Hortaliza"
